# "Generate Report for Handoff"
#
# The localization-status report is regenerated: a new handoff run
# completes for the six files that were "Ready for handoff", so:
#   - their "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#     timestamps are refreshed, and
#   - their "Priority" column is now populated with the handoff type "ht"
#     (it used to be blank while the handoff was still pending).
#
# Rows untouched: the two files that are still pending handoff
# (status column unaffected, "Latest Handoff Datetime" still the
# 0001-01-01 placeholder) keep their blank Priority.

$wb = $excel.ActiveWorkbook

$rows = @("7", "8", "10", "11", "12", "13")

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G" + $r).Value = "2016-09-01 18:24:46"
}

# --- zh-cn sheet: "Latest Handoff Datetime" (H) + "Priority" (E) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("H" + $r).Value = "2016-09-01 18:24:41"
    $zhcn.Range("E" + $r).Value = "ht"
}

# --- de-de sheet: "Latest Handoff Datetime" (H) + "Priority" (E) ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("H" + $r).Value = "2016-09-01 18:24:46"
    $dede.Range("E" + $r).Value = "ht"
}
